# Apply cryptocurrency price/volume updates to columns D and E for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    # Force the cell to remain a text value (matches source inlineStr cells)
    # even when the string looks like a number (e.g. "0.999"),
    # then restore the default "Normal" style so no stray formatting is left behind.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "62.255.72"
$ws.Range("E2").Value = "  -2.37%  "
Set-TextValue $ws "D3" "3.002.07"
$ws.Range("E3").Value = "  -2.47%  "
Set-TextValue $ws "D4" "0.999"
$ws.Range("E4").Value = "  -0.13%  "
Set-TextValue $ws "D5" "581.34"
$ws.Range("E5").Value = "  -0.96%  "
Set-TextValue $ws "D6" "146.92"
$ws.Range("E6").Value = "  -5.60%  "
$ws.Range("E7").Value = "  +0.03%  "
Set-TextValue $ws "D8" "0.520"
$ws.Range("E8").Value = "  -3.27%  "
Set-TextValue $ws "D9" "3.003.40"
$ws.Range("E9").Value = "  -2.37%  "
$ws.Range("E10").Value = "  -5.23%  "
Set-TextValue $ws "D11" "5.65"
$ws.Range("E11").Value = "  -3.40%  "
$ws.Range("E12").Value = "  -2.24%  "
Set-TextValue $ws "D13" "0.0000228"
$ws.Range("E13").Value = "  -4.06%  "
Set-TextValue $ws "D14" "34.64"
$ws.Range("E14").Value = "  -5.83%  "
$ws.Range("E15").Value = "  +1.39%  "
Set-TextValue $ws "D16" "3.493.67"
$ws.Range("E16").Value = "  -2.53%  "
$ws.Range("E17").Value = "  -2.07%  "
Set-TextValue $ws "D18" "62.210.46"
$ws.Range("E18").Value = "  -2.29%  "
Set-TextValue $ws "D19" "2.999.71"
$ws.Range("E19").Value = "  -2.51%  "
Set-TextValue $ws "D20" "454.64"
$ws.Range("E20").Value = "  -3.35%  "
Set-TextValue $ws "D21" "13.84"
$ws.Range("E21").Value = "  -3.20%  "
Set-TextValue $ws "D22" "0.674"
$ws.Range("E22").Value = "  -4.19%  "
$ws.Range("E23").Value = "  -2.93%  "
Set-TextValue $ws "D24" "2.30"
$ws.Range("E24").Value = "  -5.08%  "
Set-TextValue $ws "D25" "79.88"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("E26").Value = "  -4.40%  "
Set-TextValue $ws "D27" "10.12"
$ws.Range("E27").Value = "  -2.01%  "
Set-TextValue $ws "D28" "0.998"
$ws.Range("E28").Value = "  -0.24%  "
Set-TextValue $ws "D29" "0.999"
$ws.Range("E29").Value = "  -0.09%  "
Set-TextValue $ws "D30" "7.17"
$ws.Range("E30").Value = "  -2.89%  "
$ws.Range("E31").Value = "  -1.84%  "
$ws.Range("E32").Value = "  -1.84%  "
Set-TextValue $ws "D33" "26.90"
$ws.Range("E33").Value = "  -0.73%  "
Set-TextValue $ws "D35" "1.02"
$ws.Range("E35").Value = "  -1.71%  "
Set-TextValue $ws "D36" "0.0₃0793"
$ws.Range("E36").Value = "  -4.14%  "
Set-TextValue $ws "D37" "5.72"
$ws.Range("E37").Value = "  -4.28%  "
$ws.Range("E38").Value = "  -4.06%  "
Set-TextValue $ws "D39" "50.27"
$ws.Range("E39").Value = "  -0.66%  "
Set-TextValue $ws "D40" "8.93"
$ws.Range("E40").Value = "  -2.14%  "
Set-TextValue $ws "D41" "2.89"
$ws.Range("E41").Value = "  -11.00%  "
Set-TextValue $ws "D42" "410.43"
$ws.Range("E42").Value = "  -4.82%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("E44").Value = "  -4.43%  "
$ws.Range("E45").Value = "  -2.44%  "
Set-TextValue $ws "D46" "2.760.95"
$ws.Range("E46").Value = "  -1.76%  "
Set-TextValue $ws "D47" "37.93"
$ws.Range("E47").Value = "  -4.81%  "
Set-TextValue $ws "D48" "128.15"
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("E49").Value = "  +0.05%  "
Set-TextValue $ws "D50" "0.108"
$ws.Range("E50").Value = "  -1.31%  "
Set-TextValue $ws "D51" "23.65"
$ws.Range("E51").Value = "  -5.01%  "
